$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume figures (D: Price, E: Volume(1h)) for rows 2-51.
$updates = @(
    @{ Cell = 'D2'; Value = '29.141.12' }
    @{ Cell = 'E2'; Value = '  -2.38%  ' }
    @{ Cell = 'D3'; Value = '1.837.86' }
    @{ Cell = 'E3'; Value = '  -1.63%  ' }
    @{ Cell = 'E4'; Value = '  +0.02%  ' }
    @{ Cell = 'D5'; Value = '239.87' }
    @{ Cell = 'E5'; Value = '  -2.78%  ' }
    @{ Cell = 'D6'; Value = '0.6815' }
    @{ Cell = 'E6'; Value = '  -2.78%  ' }
    @{ Cell = 'D8'; Value = '0.2996' }
    @{ Cell = 'E8'; Value = '  -3.02%  ' }
    @{ Cell = 'D9'; Value = '0.07453' }
    @{ Cell = 'E9'; Value = '  -4.25%  ' }
    @{ Cell = 'D10'; Value = '23.19' }
    @{ Cell = 'E10'; Value = '  -3.09%  ' }
    @{ Cell = 'D11'; Value = '0.07638' }
    @{ Cell = 'E11'; Value = '  -2.71%  ' }
    @{ Cell = 'D12'; Value = '1.840.11' }
    @{ Cell = 'E12'; Value = '  -1.07%  ' }
    @{ Cell = 'D13'; Value = '5.042' }
    @{ Cell = 'E13'; Value = '  -2.90%  ' }
    @{ Cell = 'D14'; Value = '0.6800' }
    @{ Cell = 'E14'; Value = '  -2.51%  ' }
    @{ Cell = 'D15'; Value = '87.85' }
    @{ Cell = 'E15'; Value = '  -5.47%  ' }
    @{ Cell = 'D16'; Value = '6.117' }
    @{ Cell = 'E16'; Value = '  -7.97%  ' }
    @{ Cell = 'D17'; Value = '29.150.05' }
    @{ Cell = 'E17'; Value = '  -2.33%  ' }
    @{ Cell = 'D18'; Value = '0.000008188' }
    @{ Cell = 'E18'; Value = '  -2.78%  ' }
    @{ Cell = 'D19'; Value = '2.087.04' }
    @{ Cell = 'E19'; Value = '  -1.21%  ' }
    @{ Cell = 'D20'; Value = '230.65' }
    @{ Cell = 'E20'; Value = '  -5.69%  ' }
    @{ Cell = 'D21'; Value = '12.52' }
    @{ Cell = 'E21'; Value = '  -2.61%  ' }
    @{ Cell = 'D22'; Value = '1.000' }
    @{ Cell = 'E22'; Value = '  +0.03%  ' }
    @{ Cell = 'D23'; Value = '7.352' }
    @{ Cell = 'E23'; Value = '  -4.14%  ' }
    @{ Cell = 'E24'; Value = '  +0.03%  ' }
    @{ Cell = 'D25'; Value = '160.50' }
    @{ Cell = 'E25'; Value = '  +0.21%  ' }
    @{ Cell = 'D26'; Value = '0.1437' }
    @{ Cell = 'E26'; Value = '  -5.21%  ' }
    @{ Cell = 'D27'; Value = '8.694' }
    @{ Cell = 'E27'; Value = '  -3.37%  ' }
    @{ Cell = 'D28'; Value = '18.07' }
    @{ Cell = 'E28'; Value = '  -2.15%  ' }
    @{ Cell = 'D29'; Value = '1.501' }
    @{ Cell = 'E29'; Value = '  -3.08%  ' }
    @{ Cell = 'D30'; Value = '4.264' }
    @{ Cell = 'E30'; Value = '  -0.77%  ' }
    @{ Cell = 'D31'; Value = '4.138' }
    @{ Cell = 'E31'; Value = '  -2.48%  ' }
    @{ Cell = 'D32'; Value = '1.192' }
    @{ Cell = 'E32'; Value = '  -1.02%  ' }
    @{ Cell = 'D33'; Value = '0.05359' }
    @{ Cell = 'E33'; Value = '  +4.84%  ' }
    @{ Cell = 'D34'; Value = '0.7544' }
    @{ Cell = 'E34'; Value = '  -4.58%  ' }
    @{ Cell = 'D35'; Value = '1.852' }
    @{ Cell = 'E35'; Value = '  -4.49%  ' }
    @{ Cell = 'E36'; Value = '  -3.11%  ' }
    @{ Cell = 'D37'; Value = '2.689' }
    @{ Cell = 'E37'; Value = '  -0.84%  ' }
    @{ Cell = 'D38'; Value = '1.310.96' }
    @{ Cell = 'E38'; Value = '  -1.99%  ' }
    @{ Cell = 'E39'; Value = '  -3.13%  ' }
    @{ Cell = 'D40'; Value = '2.724' }
    @{ Cell = 'E40'; Value = '  -1.00%  ' }
    @{ Cell = 'D41'; Value = '0.9459' }
    @{ Cell = 'E41'; Value = '  -2.22%  ' }
    @{ Cell = 'D42'; Value = '6.003' }
    @{ Cell = 'E42'; Value = '  -1.16%  ' }
    @{ Cell = 'D43'; Value = '104.56' }
    @{ Cell = 'E43'; Value = '  -2.55%  ' }
    @{ Cell = 'D44'; Value = '0.9995' }
    @{ Cell = 'E44'; Value = '  -0.02%  ' }
    @{ Cell = 'D45'; Value = '1.990.17' }
    @{ Cell = 'E45'; Value = '  -1.13%  ' }
    @{ Cell = 'E46'; Value = '  -0.28%  ' }
    @{ Cell = 'E47'; Value = '  -3.35%  ' }
    @{ Cell = 'D48'; Value = '64.33' }
    @{ Cell = 'E48'; Value = '  -1.85%  ' }
    @{ Cell = 'D49'; Value = '9.455' }
    @{ Cell = 'E49'; Value = '  -3.85%  ' }
    @{ Cell = 'D50'; Value = '1.769' }
    @{ Cell = 'E50'; Value = '  -1.89%  ' }
    @{ Cell = 'D51'; Value = '0.07655' }
    @{ Cell = 'E51'; Value = '  +15.19%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    # Force text storage so numeric-looking strings (e.g. "0.6800", "1.000")
    # keep their exact formatting instead of being coerced to numbers.
    $rng.NumberFormat = '@'
    $rng.Value = $u.Value
}
